# Update Mean (H) and Std (I) values for specific rows on the "Overall" sheet
# as part of "feat: update preferences order evalution"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overall")

$ws.Cells.Item(26, 8).Value = 0.58336
$ws.Cells.Item(26, 9).Value = 0.02851
$ws.Cells.Item(27, 8).Value = 0.01837
$ws.Cells.Item(27, 9).Value = 0.02395
$ws.Cells.Item(28, 8).Value = 0.5910300000000001
$ws.Cells.Item(28, 9).Value = 0.02481
$ws.Cells.Item(29, 8).Value = 0.03089
$ws.Cells.Item(29, 9).Value = 0.02935
$ws.Cells.Item(30, 8).Value = 0.58356
$ws.Cells.Item(30, 9).Value = 0.02713
$ws.Cells.Item(31, 8).Value = 0.01935
$ws.Cells.Item(31, 9).Value = 0.02368
$ws.Cells.Item(32, 8).Value = 0.59414
$ws.Cells.Item(32, 9).Value = 0.02651
$ws.Cells.Item(33, 8).Value = 0.03858
$ws.Cells.Item(33, 9).Value = 0.03193
$ws.Cells.Item(34, 8).Value = 0.61552
$ws.Cells.Item(34, 9).Value = 0.01264
$ws.Cells.Item(35, 8).Value = 0.00098
$ws.Cells.Item(35, 9).Value = 0.00478
$ws.Cells.Item(36, 8).Value = 0.61578
$ws.Cells.Item(36, 9).Value = 0.01244
$ws.Cells.Item(37, 8).Value = 0.00098
$ws.Cells.Item(37, 9).Value = 0.00478
$ws.Cells.Item(38, 8).Value = 0.61558
$ws.Cells.Item(38, 9).Value = 0.01279
$ws.Cells.Item(39, 8).Value = 0.00098
$ws.Cells.Item(39, 9).Value = 0.00478
$ws.Cells.Item(40, 8).Value = 0.6158400000000001
$ws.Cells.Item(40, 9).Value = 0.01262
$ws.Cells.Item(41, 8).Value = 0.00098
$ws.Cells.Item(41, 9).Value = 0.00478
$ws.Cells.Item(66, 8).Value = 0.50317
$ws.Cells.Item(66, 9).Value = 0.03149
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(68, 8).Value = 0.54526
$ws.Cells.Item(68, 9).Value = 0.02806
$ws.Cells.Item(69, 8).Value = 0.009639999999999999
$ws.Cells.Item(69, 9).Value = 0.01797
$ws.Cells.Item(70, 8).Value = 0.5081
$ws.Cells.Item(70, 9).Value = 0.02999
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(72, 8).Value = 0.54813
$ws.Cells.Item(72, 9).Value = 0.02912
$ws.Cells.Item(73, 8).Value = 0.01454
$ws.Cells.Item(73, 9).Value = 0.01929
$ws.Cells.Item(74, 8).Value = 0.56064
$ws.Cells.Item(74, 9).Value = 0.02111
$ws.Cells.Item(75, 8).Value = 0.0029
$ws.Cells.Item(75, 9).Value = 0.007860000000000001
$ws.Cells.Item(76, 8).Value = 0.57012
$ws.Cells.Item(76, 9).Value = 0.01988
$ws.Cells.Item(77, 8).Value = 0.00388
$ws.Cells.Item(77, 9).Value = 0.00889
$ws.Cells.Item(78, 8).Value = 0.56065
$ws.Cells.Item(78, 9).Value = 0.0214
$ws.Cells.Item(79, 8).Value = 0.0029
$ws.Cells.Item(79, 9).Value = 0.007860000000000001
$ws.Cells.Item(80, 8).Value = 0.57064
$ws.Cells.Item(80, 9).Value = 0.01981
$ws.Cells.Item(81, 8).Value = 0.00388
$ws.Cells.Item(81, 9).Value = 0.00889
$ws.Cells.Item(106, 8).Value = 0.58602
$ws.Cells.Item(106, 9).Value = 0.02371
$ws.Cells.Item(107, 8).Value = 0.02019
$ws.Cells.Item(107, 9).Value = 0.02204
$ws.Cells.Item(108, 8).Value = 0.59419
$ws.Cells.Item(108, 9).Value = 0.01727
$ws.Cells.Item(109, 8).Value = 0.05593
$ws.Cells.Item(109, 9).Value = 0.03891
$ws.Cells.Item(110, 8).Value = 0.5860300000000001
$ws.Cells.Item(110, 9).Value = 0.02272
$ws.Cells.Item(111, 8).Value = 0.01923
$ws.Cells.Item(111, 9).Value = 0.02025
$ws.Cells.Item(112, 8).Value = 0.59631
$ws.Cells.Item(112, 9).Value = 0.01748
$ws.Cells.Item(113, 8).Value = 0.05974
$ws.Cells.Item(113, 9).Value = 0.04347
$ws.Cells.Item(114, 8).Value = 0.61519
$ws.Cells.Item(114, 9).Value = 0.01672
$ws.Cells.Item(115, 8).Value = 0.00479
$ws.Cells.Item(115, 9).Value = 0.00957
$ws.Cells.Item(116, 8).Value = 0.61615
$ws.Cells.Item(116, 9).Value = 0.01646
$ws.Cells.Item(117, 8).Value = 0.00574
$ws.Cells.Item(117, 9).Value = 0.01223
$ws.Cells.Item(118, 8).Value = 0.61519
$ws.Cells.Item(118, 9).Value = 0.01672
$ws.Cells.Item(119, 8).Value = 0.00479
$ws.Cells.Item(119, 9).Value = 0.00957
$ws.Cells.Item(120, 8).Value = 0.61615
$ws.Cells.Item(120, 9).Value = 0.01646
$ws.Cells.Item(121, 8).Value = 0.00574
$ws.Cells.Item(121, 9).Value = 0.01223
$ws.Cells.Item(146, 8).Value = 0.50785
$ws.Cells.Item(146, 9).Value = 0.04348
$ws.Cells.Item(147, 8).Value = 0.008710000000000001
$ws.Cells.Item(147, 9).Value = 0.02158
$ws.Cells.Item(148, 8).Value = 0.54376
$ws.Cells.Item(148, 9).Value = 0.03925
$ws.Cells.Item(149, 8).Value = 0.01738
$ws.Cells.Item(149, 9).Value = 0.02607
$ws.Cells.Item(150, 8).Value = 0.50828
$ws.Cells.Item(150, 9).Value = 0.04328
$ws.Cells.Item(151, 8).Value = 0.008710000000000001
$ws.Cells.Item(151, 9).Value = 0.02158
$ws.Cells.Item(152, 8).Value = 0.54843
$ws.Cells.Item(152, 9).Value = 0.0387
$ws.Cells.Item(153, 8).Value = 0.02794
$ws.Cells.Item(153, 9).Value = 0.02871
$ws.Cells.Item(154, 8).Value = 0.56306
$ws.Cells.Item(154, 9).Value = 0.03087
$ws.Cells.Item(155, 8).Value = 0.00581
$ws.Cells.Item(155, 9).Value = 0.0158
$ws.Cells.Item(156, 8).Value = 0.57157
$ws.Cells.Item(156, 9).Value = 0.02729
$ws.Cells.Item(157, 8).Value = 0.00678
$ws.Cells.Item(157, 9).Value = 0.01616
$ws.Cells.Item(158, 8).Value = 0.56338
$ws.Cells.Item(158, 9).Value = 0.03099
$ws.Cells.Item(159, 8).Value = 0.00581
$ws.Cells.Item(159, 9).Value = 0.0158
$ws.Cells.Item(160, 8).Value = 0.57165
$ws.Cells.Item(160, 9).Value = 0.0273
$ws.Cells.Item(161, 8).Value = 0.00678
$ws.Cells.Item(161, 9).Value = 0.01616
